$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A47").Value = "strange tiny shift like Год выпуска [___] and Год выпуска   [___] In addGoodsToWarehouse"
$ws.Range("A48").Value = "sort specific types"

$ws.Range("A47:A48").WrapText = $true

$ws.Range("A48").Select()
